$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("company_name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "not found: company_name" }
$para = $rng.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4AFD26F7" w14:textId="4BCB8DEF" w:rsidR="00F75F2B" w:rsidRPr="00310659" w:rsidRDefault="00742CDB" w:rsidP="00842BE6"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00310659"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="000000"/><w:szCs w:val="21"/></w:rPr><w:t>company_name</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para.InsertXML($xml)

$rng = $d.Content
$found = $rng.Find.Execute("bank_name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "not found: bank_name" }
$para = $rng.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="75900985" w14:textId="7654CFC0" w:rsidR="00842BE6" w:rsidRPr="00310659" w:rsidRDefault="00742CDB" w:rsidP="00F20641"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00310659"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="21"/></w:rPr><w:t>bank_name</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para.InsertXML($xml)

$rng = $d.Content
$found = $rng.Find.Execute("account_name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "not found: account_name" }
$para = $rng.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="13810A32" w14:textId="537EEC00" w:rsidR="00842BE6" w:rsidRPr="00310659" w:rsidRDefault="00742CDB" w:rsidP="00F20641"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00310659"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="21"/></w:rPr><w:t>account_name</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para.InsertXML($xml)

$rng = $d.Content
$found = $rng.Find.Execute("bank_account_no", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "not found: bank_account_no" }
$para = $rng.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="77A68E71" w14:textId="6043E4D4" w:rsidR="00842BE6" w:rsidRPr="00310659" w:rsidRDefault="00742CDB" w:rsidP="00F20641"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00310659"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="21"/></w:rPr><w:t>bank_account_no</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para.InsertXML($xml)

$rng = $d.Content
$found = $rng.Find.Execute("payment_routing_no", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "not found: payment_routing_no" }
$para = $rng.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="19553657" w14:textId="3F10B821" w:rsidR="00842BE6" w:rsidRPr="00310659" w:rsidRDefault="00742CDB" w:rsidP="00F20641"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="21"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00310659"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="21"/></w:rPr><w:t>payment_routing_no</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para.InsertXML($xml)

$rng = $d.Content
$found = $rng.Find.Execute("aba_no", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "not found: aba_no" }
$para = $rng.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="37F3A9F7" w14:textId="50C3328F" w:rsidR="00842BE6" w:rsidRPr="00310659" w:rsidRDefault="00742CDB" w:rsidP="00F20641"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00310659"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="21"/></w:rPr><w:t>aba_no</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para.InsertXML($xml)

$rng = $d.Content
$found = $rng.Find.Execute("swift_code", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "not found: swift_code" }
$para = $rng.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2BD84074" w14:textId="67057CFF" w:rsidR="00842BE6" w:rsidRPr="00310659" w:rsidRDefault="00742CDB" w:rsidP="00F20641"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00310659"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="21"/></w:rPr><w:t>swift_code</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para.InsertXML($xml)

$rng = $d.Content
$found = $rng.Find.Execute("FINAL_AMOUNT", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "not found: FINAL_AMOUNT" }
$para = $rng.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="15A08A59" w14:textId="2349D320" w:rsidR="00A122FE" w:rsidRPr="00310659" w:rsidRDefault="00A122FE" w:rsidP="00A122FE"><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r w:rsidRPr="00310659"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>TOTAL_AMOUNT</w:t></w:r></w:p>'
$para.InsertXML($xml)

Write-Output "done"